$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.823.28'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').Value = '1.633.52'
$ws.Range('E3').Value = '  +0.26%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.97'
$ws.Range('E5').Value = '  -0.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.507'
$ws.Range('E6').Value = '  -0.59%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('E8').Value = '  -0.40%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0641'
$ws.Range('E9').Value = '  +0.14%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.04'
$ws.Range('E10').Value = '  +3.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0780'
$ws.Range('E11').Value = '  +0.20%  '
$ws.Range('D12').Value = '1.666.29'
$ws.Range('E12').Value = '  +2.29%  '
$ws.Range('E13').Value = '  -0.32%  '
$ws.Range('D14').Value = '1.860.70'
$ws.Range('E14').Value = '  +0.37%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.561'
$ws.Range('E15').Value = '  +0.11%  '
$ws.Range('E16').Value = '  +1.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.06'
$ws.Range('E17').Value = '  -0.67%  '
$ws.Range('D18').Value = '25.843.85'
$ws.Range('E18').Value = '  +0.02%  '
$ws.Range('E19').Value = '  -0.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '193.69'
$ws.Range('E20').Value = '  -0.45%  '
$ws.Range('E21').Value = '  +0.82%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.90'
$ws.Range('E22').Value = '  +0.91%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.17'
$ws.Range('E23').Value = '  +2.64%  '
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('E25').Value = '  -4.63%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '138.63'
$ws.Range('E26').Value = '  -1.90%  '
$ws.Range('E27').Value = '  -4.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.81'
$ws.Range('E28').Value = '  +0.88%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.50'
$ws.Range('E29').Value = '  +0.38%  '
$ws.Range('E30').Value = '  +0.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0492'
$ws.Range('E31').Value = '  +0.68%  '
$ws.Range('E32').Value = '  -0.40%  '
$ws.Range('E33').Value = '  +1.52%  '
$ws.Range('E34').Value = '  +0.71%  '
$ws.Range('E35').Value = '  +0.55%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.901'
$ws.Range('E36').Value = '  +0.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.57'
$ws.Range('E37').Value = '  +1.31%  '
$ws.Range('D38').Value = '1.119.76'
$ws.Range('E38').Value = '  -1.16%  '
$ws.Range('E39').Value = '  -0.46%  '
$ws.Range('E40').Value = '  +0.62%  '
$ws.Range('E41').Value = '  -0.37%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.50'
$ws.Range('E42').Value = '  -1.63%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.30'
$ws.Range('E43').Value = '  +2.00%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.799'
$ws.Range('E44').Value = '  +0.38%  '
$ws.Range('D45').Value = '0.0₆0110'
$ws.Range('E45').Value = '  -1.28%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '55.35'
$ws.Range('E46').Value = '  +0.69%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.423'
$ws.Range('E47').Value = '  -4.60%  '
$ws.Range('E48').Value = '  -0.53%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.61'
$ws.Range('E49').Value = '  +0.64%  '
$ws.Range('E50').Value = '  -0.12%  '
$ws.Range('E51').Value = '  -0.20%  '
